$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Player roster shifted up by one (Nikola Jovic removed) and "Moussa Diabate"
# inserted before Walker Kessler. Rows 2-4 and 18-19 are unchanged.
$ws.Range("A5").Value = "Kel'el Ware"
$ws.Range("B5").Value = "PF,C"
$ws.Range("C5").Value = "Miami Heat"

$ws.Range("A6").Value = "Jalen Williams"
$ws.Range("B6").Value = "SG,SF,PF,C"
$ws.Range("C6").Value = "Oklahoma City Thunder"

$ws.Range("A7").Value = "Jalen Brunson"
$ws.Range("B7").Value = "PG"
$ws.Range("C7").Value = "New York Knicks"

$ws.Range("A8").Value = "Trae Young"
$ws.Range("B8").Value = "PG"
$ws.Range("C8").Value = "Atlanta Hawks"

$ws.Range("A9").Value = "Bol Bol"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Phoenix Suns"

$ws.Range("A10").Value = "Christian Braun"
$ws.Range("B10").Value = "SG,SF"
$ws.Range("C10").Value = "Denver Nuggets"

$ws.Range("A11").Value = "LeBron James"
$ws.Range("B11").Value = "SF,PF"
$ws.Range("C11").Value = "Los Angeles Lakers"

$ws.Range("A12").Value = "Kawhi Leonard"
$ws.Range("B12").Value = "SG,SF,PF"
$ws.Range("C12").Value = "LA Clippers"

$ws.Range("A13").Value = "Desmond Bane"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Memphis Grizzlies"

$ws.Range("A14").Value = "Immanuel Quickley"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Toronto Raptors"

$ws.Range("A15").Value = "Myles Turner"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Indiana Pacers"

$ws.Range("A16").Value = "Moussa Diabate"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Charlotte Hornets"

$ws.Range("A17").Value = "Walker Kessler"
$ws.Range("B17").Value = "C"
$ws.Range("C17").Value = "Utah Jazz"
